# Apply edit: reword several questions in "direct" sheet, and add a new
# "semantic" worksheet with new/rearranged questions, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Reword the questions on "direct" that were edited in place ---
$ws1.Range('A15').Value = 'What genes are involved in interactions with or regulation of Histone Deacetylase 1 (HDAC1) activity and function?'
$ws1.Range('A16').Value = 'What drugs are used to treat Autism Spectrum Disorder?'
$ws1.Range('A25').Value = 'What are the mechanisms of action of drugs that target the Histone Deacetylase 1 (HDAC1)?'
$ws1.Range('A28').Value = 'List the approval statuses of drugs that target Free Fatty Acid Receptor 2 (FFAR2).'
$ws1.Range('A31').Value = 'What are the drugs and their corresponding approval statuses associated with abortion that target the cAMP signaling pathway?'
$ws1.Range('A35').Value = 'What diseases are related to the Ras signaling pathway?'
$ws1.Range('A39').Value = 'List the biomarkers associated with malaria.'
$ws1.Range('A40').Value = 'Provide a list of drugs used in the treatment of malaria.'
$ws1.Range('A43').Value = 'What is the gene associated with the Free Fatty Acid Receptor 2 (FFAR2)?'
$ws1.Range('A44').Value = 'List all lung cancer drugs that act as inhibitors.'
$ws1.Range('A47').Value = 'List the mechanisms of action for drugs targeting the Prostaglandin F2-alpha receptor (PTGFR).'

# --- 2. Move the active-cell selection on "direct" to A5 (matches target) ---
$ws1.Range('A5').Select()

# --- 3. Give "direct" an explicit portrait page setup (matches target) ---
$ws1.PageSetup.Orientation = 1

# --- 4. Insert new worksheet "semantic" right after "direct" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = 'semantic'

# --- 5. Populate header + question rows on "semantic" ---
$ws2.Range('A1').Value = 'questions'
$ws2.Range('A2').Value = 'Which drugs are being developed to target the JAK-STAT pathway in cancer?'
$ws2.Range('A3').Value = 'List drugs that act as antagonists to the CXCR4 receptor in leukemia treatment.'
$ws2.Range('A4').Value = 'List all therapies targeting the PI3K pathway in breast cancer.'
$ws2.Range('A5').Value = 'Provide drugs that target FGFR mutations in various tumors.'
$ws2.Range('A6').Value = 'Provide drugs that target PARP.'
$ws2.Range('A7').Value = 'Which drugs target CD38 in multiple myeloma?'
$ws2.Range('A8').Value = 'Which drugs target the ALK gene fusion in lung cancer?'
$ws2.Range('A9').Value = 'Which drugs target CD20 in the treatment of B-cell lymphomas?'
$ws2.Range('A10').Value = 'Which drugs are currently under investigation as CDK4/6 inhibitors for breast cancer?'
$ws2.Range('A11').Value = 'Provide the drugs used to treat diabetes.'
$ws2.Range('A12').Value = 'List drugs that act as antagonists to IL-1 in inflammatory diseases.'
$ws2.Range('A13').Value = 'List all targets associated with diabetes, including any related pathways and relevant aspects of diabetes management.'
$ws2.Range('A14').Value = 'List all biomarkers associated with cancers.'
$ws2.Range('A15').Value = 'What major diseases are associated with the JAK-STAT pathway?'
$ws2.Range('A16').Value = 'Identify targets linked to autoimmune diseases.'
$ws2.Range('A17').Value = 'Which major pathways are associated with neurodegenerative diseases?'
$ws2.Range('A18').Value = 'What are the primary mechanisms of action for drugs used in carcinoma treatment?'
$ws2.Range('A19').Value = 'List diseases associated with tyrosine kinase activity.'
$ws2.Range('A20').Value = 'List all drugs that target adrenergic receptors.'
$ws2.Range('A21').Value = 'Give me a list of all phase 2 drugs that target metabolic pathways.'

# --- 6. Header formatting: reuse the same style as the "direct" header cell ---
$ws1.Range('A1').Copy()
$ws2.Range('A1').PasteSpecial(-4122)

# --- 7. Data-row formatting: thin border around each question cell ---
$ws2.Range('A2:A21').Borders.LineStyle = 1

# --- 8. Column width + page setup for "semantic" ---
$ws2.Columns.Item(1).ColumnWidth = 104.8
$ws2.PageSetup.Orientation = 1

# --- 9. Select A8 on "semantic" and make it the active sheet/tab ---
$ws2.Range('A8').Select()
$ws2.Activate()

